# Verizon test_data.xlsx - "Added test cases to verizon.com"
#
# Adds a new "Accessories" worksheet (after the existing "Brands" sheet)
# and populates it with the new list of bundle/accessory test assertions.

$wb = $excel.ActiveWorkbook
$brands = $wb.Worksheets.Item("Brands")

# Insert the new sheet right after "Brands"
$accessories = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $brands)
$accessories.Name = "Accessories"

# Header + data
$accessories.Range("A1").Value = "Assertions"
$accessories.Range("A2").Value = "Bundle Eligible Cases"
$accessories.Range("A3").Value = "Bundle Eligible Screen Protectors"
$accessories.Range("A4").Value = "iPad Cases"
$accessories.Range("A5").Value = "Tablet Accessories"

# Header cell uses an explicit black font color (matches the "Brands" header style)
$accessories.Range("A1").Font.Color = 0

# Size column A to fit the longest entry
$accessories.Columns.Item(1).ColumnWidth = 28.26953125

# Print setup used on the new sheet
$accessories.PageSetup.Orientation = 1

# Leave the new sheet active with the last data row selected
$accessories.Activate()
$accessories.Range("A5").Select() | Out-Null
